$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.298.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.33%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.266.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.58%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.17%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'230.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.20%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +2.81%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'63.68"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +4.56%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.447"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +11.14%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +16.86%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'56.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.73%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'26.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +19.62%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.01%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.603.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.67%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.89%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +9.58%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.837"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +5.26%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.264.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.27%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'44.052.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +5.12%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +9.40%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'73.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.28%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -2.61%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'251.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.60%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.17%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.29%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +1.70%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'3.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +25.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.18%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'172.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.95%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'20.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.01%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.74%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.39%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +3.33%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +5.73%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.42%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.77%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +7.70%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'6.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +5.93%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.08%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +3.75%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.14%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'17.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +8.34%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.86%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0963"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.17%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'97.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.93%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'TrustWalletToken"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.61%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'FTXToken"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'4.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.32%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'TerraClassic"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.000209"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -7.28%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.442.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.03%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'NEARProtocol"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'2.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.45%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Celestia"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'9.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +16.27%  "
$ws.Range("E51").Style = "Normal"
